# attention-single-hop.pptx — "Small modifications to the attention blog post"
#
# 1) The "Attention weights" textbox grows by one line to also show the
#    "p_att" label (a "p" run followed by a subscripted "att" run), and is
#    repositioned/resized to make room for it.
# 2) Two new small textboxes are added near the Query/Image vectors,
#    labelled "v" with a subscripted "Q" and "v" with a subscripted "I"
#    respectively.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) "Attention weights" textbox: reposition/resize + add "p_att" line
# ---------------------------------------------------------------------
$weights = $s.Shapes.Item("TextBox 134")

$weights.Top = 47.15503937007874
$weights.Height = 50.89220472440945

$tr = $weights.TextFrame.TextRange
$secondLine = $tr.Paragraphs(2)

# Inserting after an existing paragraph (rather than the whole TextRange)
# makes the new paragraph inherit that paragraph's exact run formatting
# (scheme color + Avenir Next), just like duplicating the "weights" line.
$null = $secondLine.InsertAfter([char]13 + "p")
$thirdLine = $tr.Paragraphs(3)
$null = $thirdLine.InsertAfter("att")

# Re-fetch the paragraph and drop the subscript ("att") onto just the
# second run, leaving "p" alone.
$thirdLine = $tr.Paragraphs(3)
$attRun = $thirdLine.Characters(2, 3)
$attRun.Font.BaselineOffset = -0.25

# ---------------------------------------------------------------------
# 2) New textbox: "v" + subscript "Q" (near the Question vector)
# ---------------------------------------------------------------------
$vq = $s.Shapes.AddTextbox(1, 28.710393700787403, 151.18023622047244, 28.67732283464567, 24.234409448818898)
$vq.Name = "TextBox 98"
$vq.Fill.Visible = 0
$vq.TextFrame.WordWrap = 0
$vq.TextFrame.AutoSize = 1

$vqTr = $vq.TextFrame.TextRange
$vqTr.Text = "vQ"
$vqTr.Font.Size = 14
$vqTr.Font.Name = "Monaco"
$vqTr.Font.NameFarEast = "Helvetica Neue"
$vqTr.Font.NameComplexScript = "Helvetica Neue"

$vqSub = $vqTr.Characters(2, 1)
$vqSub.Font.BaselineOffset = -0.25

# ---------------------------------------------------------------------
# 3) New textbox: "v" + subscript "I" (near the Image vector)
# ---------------------------------------------------------------------
$vi = $s.Shapes.AddTextbox(1, 35.644488188976375, 54.80976377952756, 28.67732283464567, 24.234409448818898)
$vi.Name = "TextBox 99"
$vi.Fill.Visible = 0
$vi.TextFrame.WordWrap = 0
$vi.TextFrame.AutoSize = 1

$viTr = $vi.TextFrame.TextRange
$viTr.Text = "vI"
$viTr.Font.Size = 14
$viTr.Font.Name = "Monaco"
$viTr.Font.NameFarEast = "Helvetica Neue"
$viTr.Font.NameComplexScript = "Helvetica Neue"

$viSub = $viTr.Characters(2, 1)
$viSub.Font.BaselineOffset = -0.25
